$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay text (match source formatting exactly)
$textCells = @("D5", "D6", "D9", "D20", "D21", "D24", "D28", "D30", "D31", "D34", "D37", "D40", "D41", "D42", "D43", "D47", "D49", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated coin data
$ws.Range("D2").Value = "64.119.84"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "3.155.77"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "604.17"
$ws.Range("E5").Value = "  -2.06%  "

$ws.Range("D6").Value = "143.95"
$ws.Range("E6").Value = "  -2.75%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.145.02"
$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("E10").Value = "  -1.69%  "

$ws.Range("E11").Value = "  -2.10%  "

$ws.Range("E12").Value = "  -1.70%  "

$ws.Range("E13").Value = "  -2.20%  "

$ws.Range("E14").Value = "  -2.53%  "

$ws.Range("D15").Value = "3.672.59"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("D17").Value = "64.161.63"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").Value = "3.152.57"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("D20").Value = "492.07"
$ws.Range("E20").Value = "  +2.29%  "

$ws.Range("D21").Value = "14.73"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("E22").Value = "  -1.98%  "

$ws.Range("E23").Value = "  -4.19%  "

$ws.Range("D24").Value = "88.11"
$ws.Range("E24").Value = "  +4.10%  "

$ws.Range("E25").Value = "  -3.56%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("E27").Value = "  -2.46%  "

$ws.Range("D28").Value = "8.23"
$ws.Range("E28").Value = "  -4.42%  "

$ws.Range("E29").Value = "  +0.54%  "

$ws.Range("D30").Value = "2.05"
$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("D31").Value = "27.68"
$ws.Range("E31").Value = "  +4.05%  "

$ws.Range("E32").Value = "  -5.06%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "2.67"
$ws.Range("E34").Value = "  -1.99%  "

$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("D37").Value = "52.69"
$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("D38").Value = "0.0₃0743"
$ws.Range("E38").Value = "  -4.78%  "

$ws.Range("E39").Value = "  -7.96%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0398"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "435.35"
$ws.Range("E41").Value = "  -5.70%  "

$ws.Range("D42").Value = "0.120"
$ws.Range("E42").Value = "  -1.09%  "

$ws.Range("D43").Value = "8.30"
$ws.Range("E43").Value = "  -1.58%  "

$ws.Range("D44").Value = "2.945.11"
$ws.Range("E44").Value = "  +3.29%  "

$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("E46").Value = "  -5.92%  "

$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  -2.06%  "

$ws.Range("D49").Value = "25.92"
$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("D51").Value = "120.29"
$ws.Range("E51").Value = "  -0.19%  "
